# Hjemme passive updated meanEMG legmaxROM
# Update the B:E columns for rows 1-3 on the active sheet ("Ark1") with the
# refreshed values, then narrow the saved selection to the edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 19.484201272969916
$ws.Range("C2").Value = 5.1676634364079881
$ws.Range("D2").Value = 6.1620420962219775
$ws.Range("E2").Value = 1.783092086562087

$ws.Range("B3").Value = 33.143232035472344
$ws.Range("C3").Value = 4.2882762550519846
$ws.Range("D3").Value = -4.9049836629851153
$ws.Range("E3").Value = 5.8290072667893034

$ws.Range("B1:E3").Select()
